# Daily cryptos.xlsx price/volume refresh (GitHub Actions bot)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.153.34'
$ws.Range('E2').Value = '  -1.43%  '
$ws.Range('D3').Value = '2.767.29'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'353.53"
$ws.Range('E5').Value = '  -0.47%  '
$ws.Range('D6').Value = "'107.38"
$ws.Range('E6').Value = '  -0.92%  '
$ws.Range('D7').Value = "'0.547"
$ws.Range('E7').Value = '  -2.52%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = "'0.582"
$ws.Range('E9').Value = '  -0.90%  '
$ws.Range('D10').Value = "'39.42"
$ws.Range('E10').Value = '  -1.47%  '
$ws.Range('E11').Value = '  +3.45%  '
$ws.Range('D12').Value = "'0.0831"
$ws.Range('E12').Value = '  -2.07%  '
$ws.Range('D13').Value = "'19.96"
$ws.Range('E13').Value = '  +3.54%  '
$ws.Range('D14').Value = "'7.51"
$ws.Range('E14').Value = '  -0.66%  '
$ws.Range('D15').Value = '3.202.40'
$ws.Range('E15').Value = '  -0.18%  '
$ws.Range('D16').Value = '2.772.15'
$ws.Range('E16').Value = '  -0.13%  '
$ws.Range('D17').Value = "'0.927"
$ws.Range('E17').Value = '  +0.09%  '
$ws.Range('D18').Value = '51.167.36'
$ws.Range('E18').Value = '  -1.23%  '
$ws.Range('D19').Value = "'7.64"
$ws.Range('E19').Value = '  +3.94%  '
$ws.Range('E20').Value = '  -0.95%  '
$ws.Range('D21').Value = "'13.06"
$ws.Range('E21').Value = '  +0.72%  '
$ws.Range('D22').Value = '0.0₃0960'
$ws.Range('E22').Value = '  -1.20%  '
$ws.Range('D23').Value = "'69.62"
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('D24').Value = "'265.32"
$ws.Range('E24').Value = '  -3.14%  '
$ws.Range('D25').Value = "'2.72"
$ws.Range('E25').Value = '  -0.09%  '
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('D27').Value = "'25.91"
$ws.Range('E27').Value = '  -2.01%  '
$ws.Range('D28').Value = "'0.161"
$ws.Range('E28').Value = '  +11.89%  '
$ws.Range('D29').Value = "'10.14"
$ws.Range('E29').Value = '  +0.58%  '
$ws.Range('D30').Value = "'2.21"
$ws.Range('E30').Value = '  -0.85%  '
$ws.Range('D31').Value = "'34.75"
$ws.Range('E31').Value = '  +3.40%  '
$ws.Range('D32').Value = "'51.74"
$ws.Range('E32').Value = '  +0.87%  '
$ws.Range('D33').Value = "'6.08"
$ws.Range('E33').Value = '  +7.01%  '
$ws.Range('B34').Value = 'RenderToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D34').Value = "'5.52"
$ws.Range('E34').Value = '  +3.99%  '
$ws.Range('B35').Value = 'VeChain'
$ws.Range('C35').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D35').Value = "'0.0440"
$ws.Range('E35').Value = '  -4.51%  '
$ws.Range('D36').Value = "'0.0828"
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('D38').Value = "'18.13"
$ws.Range('E38').Value = '  +0.30%  '
$ws.Range('E39').Value = '  -1.45%  '
$ws.Range('E40').Value = '  -1.53%  '
$ws.Range('E41').Value = '  -0.76%  '
$ws.Range('E42').Value = '  -0.48%  '
$ws.Range('D43').Value = "'120.73"
$ws.Range('E43').Value = '  -1.43%  '
$ws.Range('D44').Value = "'22.03"
$ws.Range('E44').Value = '  +1.78%  '
$ws.Range('D45').Value = "'2.18"
$ws.Range('E45').Value = '  -2.62%  '
$ws.Range('D46').Value = '2.086.32'
$ws.Range('E46').Value = '  +1.59%  '
$ws.Range('D47').Value = "'3.23"
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('E48').Value = '  -0.64%  '
$ws.Range('B49').Value = 'SEI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range('D49').Value = "'0.912"
$ws.Range('E49').Value = '  -0.37%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').Value = "'5.42"
$ws.Range('E50').Value = '  -4.54%  '
$ws.Range('D51').Value = "'1.30"
$ws.Range('E51').Value = '  +7.13%  '
